# SM added low stock
# Applies three edits inside the "text files details" table:
#  1. Adds a new highlighted paragraph "boolean stockLevel; (5 May)" right
#     after the "double price;" line in the "items" row.
#  2. Splits the "suppliers" row-label run into "s" + "uppliers" (with a
#     stray spell-check proofErr marker between them).
#  3. Highlights (yellow) the struck-through "String itemId;" remark in the
#     "suppliers" row, and adds an eastAsia font hint to its paragraph mark.

$d = $word.ActiveDocument

$pkgOpen = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------
# 1) New paragraph after "double price;" (paragraph 7 in the items row)
# ---------------------------------------------------------------------
$pricePara = $d.Paragraphs(7)
if (-not $pricePara.Range.Text.StartsWith("double price;")) {
    throw "Unexpected paragraph 7 text: $($pricePara.Range.Text)"
}
$pricePara.Range.InsertParagraphAfter()

$stockPara = $d.Paragraphs(8)
$stockXml = $pkgOpen + `
  '<w:p>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>boolean</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>stockLevel</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>;</w:t></w:r>' + `
    '<w:r><w:t xml:space="preserve"> (5 May)</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
  '</w:p>' + `
  $pkgClose
$stockPara.Range.InsertXML($stockXml)

# ---------------------------------------------------------------------
# 2) Split "suppliers" row-label run into "s" / "uppliers"
# ---------------------------------------------------------------------
$supPara = $d.Paragraphs(10)
if (-not $supPara.Range.Text.StartsWith("suppliers")) {
    throw "Unexpected paragraph 10 text: $($supPara.Range.Text)"
}
$supXml = $pkgOpen + `
  '<w:p>' + `
    '<w:r><w:t>s</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>uppliers</w:t></w:r>' + `
  '</w:p>' + `
  $pkgClose
$supPara.Range.InsertXML($supXml)

# ---------------------------------------------------------------------
# 3) Highlight the struck-through "String itemId;" remark paragraph
# ---------------------------------------------------------------------
$itemIdPara = $d.Paragraphs(13)
if (-not $itemIdPara.Range.Text.StartsWith("String itemId; no need Item ID")) {
    throw "Unexpected paragraph 13 text: $($itemIdPara.Range.Text)"
}
$itemIdXml = $pkgOpen + `
  '<w:p>' + `
    '<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:strike/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">String </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>itemId</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t>;</w:t></w:r>' + `
    '<w:r><w:rPr><w:strike/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>no need Item ID as it might repeat the supplier</w:t></w:r>' + `
    '<w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
  '</w:p>' + `
  $pkgClose
$itemIdPara.Range.InsertXML($itemIdXml)

Write-Output "done"
